$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148; existing rows 148:196 shift down to 149:197.
$ws.Rows("148").Insert()

# Populate the newly inserted row 148 with the new record.
$ws.Range("A148").Value = 5
$ws.Range("B148").Value = "Macroferia Regional de Talca"
$ws.Range("C148").Value = "Maule"
$ws.Range("D148").Value = 44468
$ws.Range("E148").Value = 7
$ws.Range("F148").Value = 100112023
$ws.Range("G148").Value = "Brócoli"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 3000
$ws.Range("K148").Value = 600
$ws.Range("L148").Value = 600
$ws.Range("M148").Value = 600
$ws.Range("N148").Value = "`$/unidad"
$ws.Range("O148").Value = "Región del Maule"
$ws.Range("P148").Value = 600
$ws.Range("Q148").Value = 1
$ws.Range("R148").Value = "Hortaliza"
